# Append a new data row (row 52) to each of the four sheets, mirroring the
# structure of the existing rows (time stamp, raw hex fields, and their
# decoded numeric counterparts).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: DE_LFT_#1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A52").Value = 45838.43719907408
$ws1.Range("A52").NumberFormat = $ws1.Range("A51").NumberFormat
$ws1.Range("B52").Value = "0x01,0x7c"
$ws1.Range("C52").Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws1.Range("D52").Value = "0x01,0x60"
$ws1.Range("E52").Value = "0x14"
$ws1.Range("F52").Value = 380
$ws1.Range("G52").Value = [double]"7.598631275147109e+23"
$ws1.Range("H52").Value = 352
$ws1.Range("I52").Value = 14

# --- Sheet 2: DE_LFT_#2 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A52").Value = 45838.43719907408
$ws2.Range("A52").NumberFormat = $ws2.Range("A51").NumberFormat
$ws2.Range("B52").Value = "0x01,0x7c"
$ws2.Range("C52").Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws2.Range("D52").Value = "0x01,0x64"
$ws2.Range("E52").Value = "0xe"
$ws2.Range("F52").Value = 380
$ws2.Range("G52").Value = [double]"5.68432987514711e+23"
$ws2.Range("H52").Value = 356
$ws2.Range("I52").Value = 14

# --- Sheet 3: DE_PLT_#1 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A52").Value = 45838.43719907408
$ws3.Range("A52").NumberFormat = $ws3.Range("A51").NumberFormat
$ws3.Range("B52").Value = "0x00,0x82"
$ws3.Range("C52").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Range("D52").Value = "0x00,0x7E"
$ws3.Range("E52").Value = "0x7"
$ws3.Range("F52").Value = 130
$ws3.Range("G52").Value = [double]"5.68631262647114e+23"
$ws3.Range("H52").Value = 126
$ws3.Range("I52").Value = 7

# --- Sheet 4: DE_PLT_#2 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A52").Value = 45838.43719907408
$ws4.Range("A52").NumberFormat = $ws4.Range("A51").NumberFormat
$ws4.Range("B52").Value = "0x00,0x82"
$ws4.Range("C52").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Range("D52").Value = "0x00,0x7E"
$ws4.Range("E52").Value = "0x3"
$ws4.Range("F52").Value = 130
$ws4.Range("G52").Value = [double]"9.85046333984776e+23"
$ws4.Range("H52").Value = 126
$ws4.Range("I52").Value = 3
